$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Range("A2").Value = "FPK12School65086"
$ws.Range("B2").Value = "FPK12Classroom67383"
$ws.Range("C2").Value = "FPK12Section65747"
$ws.Range("E3").Value = "26235"
$ws.Range("E4").Value = "58623"
$ws.Range("E5").Value = "13127"
